$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell='D2'; Value='42.042.85'; ForceText=$false},
    @{Cell='E2'; Value='  -0.21%  '; ForceText=$false},
    @{Cell='D3'; Value='2.258.14'; ForceText=$false},
    @{Cell='E3'; Value='  -0.43%  '; ForceText=$false},
    @{Cell='D4'; Value='1.00'; ForceText=$true},
    @{Cell='E4'; Value='  +0.04%  '; ForceText=$false},
    @{Cell='D5'; Value='304.85'; ForceText=$true},
    @{Cell='E5'; Value='  -0.37%  '; ForceText=$false},
    @{Cell='D6'; Value='96.01'; ForceText=$true},
    @{Cell='E6'; Value='  +2.07%  '; ForceText=$false},
    @{Cell='E7'; Value='  -1.12%  '; ForceText=$false},
    @{Cell='E8'; Value='  +0.01%  '; ForceText=$false},
    @{Cell='D9'; Value='0.490'; ForceText=$true},
    @{Cell='E9'; Value='  +0.77%  '; ForceText=$false},
    @{Cell='E10'; Value='  +5.06%  '; ForceText=$false},
    @{Cell='E11'; Value='  -1.68%  '; ForceText=$false},
    @{Cell='E12'; Value='  -0.18%  '; ForceText=$false},
    @{Cell='E13'; Value='  +2.10%  '; ForceText=$false},
    @{Cell='D14'; Value='2.608.15'; ForceText=$false},
    @{Cell='D15'; Value='14.47'; ForceText=$true},
    @{Cell='E15'; Value='  +0.81%  '; ForceText=$false},
    @{Cell='D16'; Value='2.258.52'; ForceText=$false},
    @{Cell='E16'; Value='  -0.19%  '; ForceText=$false},
    @{Cell='D17'; Value='0.791'; ForceText=$true},
    @{Cell='E17'; Value='  +0.36%  '; ForceText=$false},
    @{Cell='D18'; Value='41.926.54'; ForceText=$false},
    @{Cell='E18'; Value='  -0.20%  '; ForceText=$false},
    @{Cell='D19'; Value='12.34'; ForceText=$true},
    @{Cell='E19'; Value='  -2.95%  '; ForceText=$false},
    @{Cell='E20'; Value='  -1.47%  '; ForceText=$false},
    @{Cell='E21'; Value='  -0.72%  '; ForceText=$false},
    @{Cell='D22'; Value='68.45'; ForceText=$true},
    @{Cell='E22'; Value='  +0.51%  '; ForceText=$false},
    @{Cell='D23'; Value='237.07'; ForceText=$true},
    @{Cell='E24'; Value='  -2.27%  '; ForceText=$false},
    @{Cell='E25'; Value='  +0.06%  '; ForceText=$false},
    @{Cell='E26'; Value='  -1.42%  '; ForceText=$false},
    @{Cell='D27'; Value='23.51'; ForceText=$true},
    @{Cell='E27'; Value='  -2.05%  '; ForceText=$false},
    @{Cell='D28'; Value='36.43'; ForceText=$true},
    @{Cell='E28'; Value='  +4.02%  '; ForceText=$false},
    @{Cell='E29'; Value='  +1.68%  '; ForceText=$false},
    @{Cell='D30'; Value='9.43'; ForceText=$true},
    @{Cell='E30'; Value='  -2.62%  '; ForceText=$false},
    @{Cell='D31'; Value='160.84'; ForceText=$true},
    @{Cell='E31'; Value='  +0.13%  '; ForceText=$false},
    @{Cell='B32'; Value='FirstDigitalUSD'; ForceText=$false},
    @{Cell='C32'; Value='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; ForceText=$false},
    @{Cell='D32'; Value='1.00'; ForceText=$true},
    @{Cell='E32'; Value='  +0.08%  '; ForceText=$false},
    @{Cell='B33'; Value='Filecoin'; ForceText=$false},
    @{Cell='C33'; Value='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; ForceText=$false},
    @{Cell='D33'; Value='5.20'; ForceText=$true},
    @{Cell='E33'; Value='  -2.51%  '; ForceText=$false},
    @{Cell='D34'; Value='3.16'; ForceText=$true},
    @{Cell='E34'; Value='  +2.78%  '; ForceText=$false},
    @{Cell='E35'; Value='  -1.37%  '; ForceText=$false},
    @{Cell='D36'; Value='17.00'; ForceText=$true},
    @{Cell='E36'; Value='  -0.27%  '; ForceText=$false},
    @{Cell='E37'; Value='  +0.08%  '; ForceText=$false},
    @{Cell='E38'; Value='  -1.67%  '; ForceText=$false},
    @{Cell='E39'; Value='  -2.07%  '; ForceText=$false},
    @{Cell='E40'; Value='  +0.20%  '; ForceText=$false},
    @{Cell='D41'; Value='3.98'; ForceText=$true},
    @{Cell='E41'; Value='  -2.19%  '; ForceText=$false},
    @{Cell='E42'; Value='  +0.68%  '; ForceText=$false},
    @{Cell='D43'; Value='1.959.98'; ForceText=$false},
    @{Cell='E43'; Value='  -2.14%  '; ForceText=$false},
    @{Cell='D44'; Value='0.0280'; ForceText=$true},
    @{Cell='E44'; Value='  -0.94%  '; ForceText=$false},
    @{Cell='D45'; Value='18.53'; ForceText=$true},
    @{Cell='E45'; Value='  -6.16%  '; ForceText=$false},
    @{Cell='D46'; Value='9.91'; ForceText=$true},
    @{Cell='E46'; Value='  -2.68%  '; ForceText=$false},
    @{Cell='E47'; Value='  -0.67%  '; ForceText=$false},
    @{Cell='D48'; Value='53.25'; ForceText=$true},
    @{Cell='E48'; Value='  -0.43%  '; ForceText=$false},
    @{Cell='D49'; Value='72.29'; ForceText=$true},
    @{Cell='E49'; Value='  -1.19%  '; ForceText=$false},
    @{Cell='D50'; Value='91.14'; ForceText=$true},
    @{Cell='E50'; Value='  -0.42%  '; ForceText=$false},
    @{Cell='E51'; Value='  -2.38%  '; ForceText=$false}
)

foreach ($chg in $changes) {
    $rng = $ws.Range($chg.Cell)
    if ($chg.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $chg.Value
}

Write-Host "Applied $($changes.Count) cell updates"